$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Price($cellRef, $text) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
}

function Set-Text($cellRef, $text) {
    $ws.Range($cellRef).Value = $text
}

# --- Simple price-only (column D) updates ---
Set-Price "D2" "247.61"
Set-Price "D3" "21.78"
Set-Price "D4" "5.423"
Set-Price "D5" "0.05691"
Set-Price "D6" "3.384"
Set-Price "D7" "0.8087"
Set-Price "D8" "1.022"
Set-Price "D9" "0.1452"
Set-Price "D10" "0.07516"
Set-Price "D11" "0.03160"
Set-Price "D12" "0.03051"
Set-Price "D13" "0.09259"

# --- Row 14: price + volume(1h) text tweak ---
Set-Price "D14" "3.624"
Set-Text "E14" "13MCDexMCB"

Set-Price "D15" "0.001645"
Set-Price "D16" "0.04704"

# --- Rows 17-25: coin list rotated up by one slot ---
Set-Text  "B17" "TigerCash"
Set-Text  "C17" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-Price "D17" "0.006351"
Set-Text  "E17" "16TigerCashTCH"

Set-Text  "B18" "HotbitToken"
Set-Text  "C18" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-Price "D18" "0.005029"
Set-Text  "E18" "17HotbitTokenHTB"

Set-Text  "B19" "BitKan"
Set-Text  "C19" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-Price "D19" "0.001043"
Set-Text  "E19" "18BitKanKAN"

Set-Text  "B20" "NitroEx"
Set-Text  "C20" "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-Price "D20" "0.0001501"
Set-Text  "E20" "19NitroExNTX"

Set-Text  "B21" "UpBots"
Set-Text  "C21" "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-Price "D21" "0.0003102"
Set-Text  "E21" "20UpBotsUBXT"

Set-Text  "B22" "LEO"
Set-Text  "C22" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-Price "D22" "3.769"
Set-Text  "E22" "21LEOLEO"

Set-Text  "B23" "KuCoinToken"
Set-Text  "C23" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-Price "D23" "6.408"
Set-Text  "E23" "22KuCoinTokenKCS"

Set-Text  "B24" "BTSEToken"
Set-Text  "C24" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-Price "D24" "2.099"
Set-Text  "E24" "23BTSETokenBTSE"

Set-Text  "B25" "One"
Set-Text  "C25" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-Price "D25" "0.0005863"
Set-Text  "E25" "24OneONE"

Set-Price "D26" "0.3324"
Set-Price "D27" "0.1306"

Set-Price "D40" "0.04051"

# --- Row 41: price + volume(1h) text tweak ---
Set-Price "D41" "0.006971"
Set-Text  "E41" "40KickTokenKICKBestin24h"

# --- Rows 42-43: BKEXToken / CEJI swap ---
Set-Text  "B42" "CEJI"
Set-Text  "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-Price "D42" "0.003502"
Set-Text  "E42" "41CEJICEJI"

Set-Text  "B43" "BKEXToken"
Set-Text  "C43" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-Price "D43" "0.1042"
Set-Text  "E43" "42BKEXTokenBKK"

Set-Price "D44" "0.008517"
Set-Price "D45" "0.00005939"
Set-Price "D46" "0.00000000750"

# --- Row 47: price + volume(1h) text tweak ---
Set-Price "D47" "0.0005503"
Set-Text  "E47" "46ACDXExchangeACXTWorstin24h"

Set-Price "D48" "0.6828"
Set-Price "D49" "0.006978"
Set-Price "D50" "0.00002101"
